$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.71"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.13"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.417"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05897"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.433"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.518"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8077"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9337"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1417"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07342"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03325"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03064"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09352"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.859"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001582"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04685"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005898"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005986"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001261"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004905"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006798"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.570"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.143"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3234"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1294"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03971"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006188"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1071"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002770"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009475"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005205"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6698"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002330"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
